$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The player roster table (A: Player, B: Position, C: Team) was re-ordered.
# Write the final values for each data row (rows 2-18) directly.

$data = @(
    @("Tyler Herro",    "PG,SG",        "Miami Heat"),
    @("Mikal Bridges",  "SG,SF,PF",     "New York Knicks"),
    @("Shaedon Sharpe", "SG,SF",        "Portland Trail Blazers"),
    @("Scottie Barnes", "PG,SG,SF,PF",  "Toronto Raptors"),
    @("Aaron Gordon",   "PF,C",         "Denver Nuggets"),
    @("Brook Lopez",    "C",            "Milwaukee Bucks"),
    @("Ja Morant",      "PG",           "Memphis Grizzlies"),
    @("Nick Richards",  "C",            "Phoenix Suns"),
    @("Kevon Looney",   "PF,C",         "Golden State Warriors"),
    @("Dillon Brooks",  "SG,SF",        "Houston Rockets"),
    @("De'Aaron Fox",   "PG",           "Sacramento Kings"),
    @("Josh Giddey",    "PG,SG,SF",     "Chicago Bulls"),
    @("Miles Bridges",  "SF,PF",        "Charlotte Hornets"),
    @("Nikola Vucevic", "PF,C",         "Chicago Bulls"),
    @("DeMar DeRozan",  "SF,PF",        "Sacramento Kings"),
    @("Luka Doncic",    "PG,SG",        "Dallas Mavericks"),
    @("Evan Mobley",    "PF,C",         "Cleveland Cavaliers")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
